$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use the existing C30 cell (already date-formatted) as a style template
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null

$ws.Range("B28").Value = "Complete"
$ws.Range("B28").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("B29").Value = "Complete"
$ws.Range("B29").Font.Color = $ws.Range("B30").Font.Color

$ws.Range("C28").Value = 40853
$ws.Range("C29").Value = 40853
$ws.Range("C31").Value = 40850

$ws.Range("C31").Select()
